# Cambios de Exportación y modificaciones de nombres de partidas
#
# - Updates the "partida" total in M1
# - Refreshes the export/update timestamp stored in P1
# - Moves the color value that used to live in U1/V1 (row 1) so that it now
#   only lives in V1 (U1 is cleared) with the updated color value
# - Removes the second "partida" row (row 2) entirely, which also drops the
#   now-unused shared strings (the HTML block, "#48abe6", "Puertas" and its
#   two timestamps)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the numeric total for the remaining item (row 1)
$ws.Range("M1").Value = 1307.236

# Refresh the timestamp held in P1
$ws.Range("P1").Value = "2018-09-11 18:52:04"

# U1 (old color value) is no longer used; the color now only lives in V1
$ws.Range("U1").ClearContents()
$ws.Range("V1").Value = "#ffffff"

# Drop the second item entirely (row 2)
$ws.Rows("2:2").Delete()
